$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 3.9
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 1.91
$ws.Range("L4").Value = 2.6
$ws.Range("AA4").Value = 29
$ws.Range("AB4").Value = 34
$ws.Range("G6").Value = 3.25
$ws.Range("H6").Value = 3.4
$ws.Range("J6").Value = 3.6
$ws.Range("K6").Value = 2.3
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 4
$ws.Range("Q6").Value = 1.75
$ws.Range("R6").Value = 2.05
$ws.Range("S6").Value = 1.33
$ws.Range("T6").Value = 3.25
$ws.Range("U6").Value = 1.57
$ws.Range("V6").Value = 2.25
$ws.Range("AB6").Value = 29
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 6.5
$ws.Range("AF6").Value = 41
$ws.Range("AH6").Value = 9.5
$ws.Range("AP6").Value = 23
$ws.Range("AR6").Value = 67
$ws.Range("AT6").Value = 3.25
$ws.Range("AU6").Value = 7.5
$ws.Range("AW6").Value = 4.33
$ws.Range("AY6").Value = 19
$ws.Range("AZ6").Value = 41
$ws.Range("BC6").Value = 451
$ws.Range("Q7").Value = 1.88
$ws.Range("R7").Value = 1.93
$ws.Range("W7").Value = 6
$ws.Range("Y7").Value = 10
$ws.Range("AH7").Value = 17
$ws.Range("AJ7").Value = 26
$ws.Range("AK7").Value = 126
$ws.Range("AR7").Value = 41
$ws.Range("AZ7").Value = 301
$ws.Range("G8").Value = 2.4
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 2.9
$ws.Range("K8").Value = 1.95
$ws.Range("Q8").Value = 2.4
$ws.Range("R8").Value = 1.53
$ws.Range("S8").Value = 1.53
$ws.Range("T8").Value = 2.38
$ws.Range("W8").Value = 6.5
$ws.Range("AB8").Value = 41
$ws.Range("AD8").Value = 6
$ws.Range("AE8").Value = 19
$ws.Range("AJ8").Value = 12
$ws.Range("AL8").Value = 29
$ws.Range("AT8").Value = 2.38
$ws.Range("AU8").Value = 9
$ws.Range("AX8").Value = 19
$ws.Range("BA8").Value = 101
$ws.Range("G10").Value = 1.4
$ws.Range("H10").Value = 4.33
$ws.Range("I10").Value = 8.5
$ws.Range("J10").Value = 1.95
$ws.Range("L10").Value = 8.5
$ws.Range("W10").Value = 5
$ws.Range("Z10").Value = 8.5
$ws.Range("AE10").Value = 29
$ws.Range("AJ10").Value = 26
$ws.Range("AM10").Value = 81
$ws.Range("AN10").Value = 3.1
$ws.Range("AQ10").Value = 21
$ws.Range("AS10").Value = 251
$ws.Range("AW10").Value = 9
$ws.Range("AZ10").Value = 251
$ws.Range("BA10").Value = 301
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("G12").Value = 3.7
$ws.Range("J12").Value = 3.95
$ws.Range("L12").Value = 2.47
$ws.Range("U12").Value = 1.6
$ws.Range("V12").Value = 2.07
$ws.Range("W12").Value = 13.5
$ws.Range("Y12").Value = 12
$ws.Range("AA12").Value = 30
$ws.Range("AB12").Value = 30
$ws.Range("AC12").Value = 12
$ws.Range("AD12").Value = 6.9
$ws.Range("AH12").Value = 8.25
$ws.Range("AI12").Value = 9.75
$ws.Range("AO12").Value = 19.5
$ws.Range("AP12").Value = 23
$ws.Range("AQ12").Value = 90
$ws.Range("AR12").Value = 110
$ws.Range("AT12").Value = 2.92
$ws.Range("AX12").Value = 9.5
$ws.Range("AY12").Value = 16.5
$ws.Range("AZ12").Value = 35
$ws.Range("BA12").Value = 60
$ws.Range("J19").Value = 2.1
$ws.Range("K19").Value = 2.12
$ws.Range("M19").Value = 1.02
$ws.Range("S19").Value = 1.4
$ws.Range("T19").Value = 2.52
$ws.Range("Y19").Value = 8
$ws.Range("AB19").Value = 28
$ws.Range("AO19").Value = 7.4
$ws.Range("AP19").Value = 18.5
$ws.Range("AQ19").Value = 24
$ws.Range("AS19").Value = 300
$ws.Range("AT19").Value = 2.5
$ws.Range("AU19").Value = 8
$ws.Range("AY19").Value = 40
